$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$values = @(
    @("68÷6=11, 2", "95÷5=19, 0", "42÷7=6, 0", "75÷2=37, 1", "83÷8=10, 3"),
    @("35÷4=8, 3", "41÷7=5, 6", "60÷5=12, 0", "31÷8=3, 7", "21÷4=5, 1"),
    @("16÷3=5, 1", "29÷5=5, 4", "98÷7=14, 0", "72÷9=8, 0", "84÷7=12, 0"),
    @("62÷8=7, 6", "74÷5=14, 4", "60÷7=8, 4", "26÷4=6, 2", "93÷7=13, 2"),
    @("67÷2=33, 1", "90÷6=15, 0", "51÷9=5, 6", "78÷9=8, 6", "15÷5=3, 0")
)

$rows = @(1, 5, 9, 13, 17)

for ($i = 0; $i -lt 5; $i++) {
    $rowIndex = $rows[$i]
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($rowIndex, $c)
        $rng = $cell.Range
        $rng.End = $rng.End - 1
        $rng.Text = $values[$i][$c - 1]
    }
}
